$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF2").Value = 1.25

$ws.Range("AE3").Value = 1.1000000000000001
$ws.Range("AF3").Value = 1
$ws.Range("AH3").Value = 0.75

$ws.Range("AE4").Value = 1.5
$ws.Range("AF4").Value = 1.25
$ws.Range("AG4").Value = 1
$ws.Range("AH4").Value = 0.75

$ws.Range("AE5").Value = 1.3
$ws.Range("AF5").Value = 1
$ws.Range("AG5").Value = 0.8
$ws.Range("AH5").Value = 0.75

$ws.Range("I14").Select()
